$wb = $excel.ActiveWorkbook

# --- Add the new "債務" (debt) sheet as the last tab ---
$template = $wb.Worksheets.Item("保險")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "債務"

# Keep the "date" column (J) from being auto-parsed into a date serial by
# pre-formatting it as text; the style gets normalised back to the template
# below once the literal text value is safely stored.
$ws.Range("J2:J8").NumberFormat = "@"

# --- Header row values (column order matches the source workbook) ---
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data, written column-major (top-to-bottom within each column, then
# moving to the next column) so freshly-minted shared strings land in the
# same slots as the original export. ---
$species  = @("房屋貸款","貸款","貸款","貸款","貸款","貸款","貸款")
$debtor   = @("黃昭順","王崇儀","王崇儀","王崇儀","王崇儀","王崇儀","王崇儀")
$owner    = @("大商業眾銀行新生分行臺北市中正區忠孝","幸福人壽保險股份有限公司臺北市中正區","幸福人壽保險股份有限公司臺北市中正區","幸福人壽保險股份有限公司臺北市中正區","幸福人壽保險股份有限公司臺北市中正區","幸福人壽保險股份有限公司臺北市中正區","幸福人壽保險股份有限公司臺北市中正區")
$total    = @(1119648,600000,1060000,500405,800820,500131,300000)
$regdate  = @("94年01月06日","95年12月25日","96年04月12曰","96年05月28日","96年12月11曰","100年12月08日","100年12月30日")
$regreas  = @("購屋貸款","般借款","般借款","般借款","般借款","一般借款","般借款")
$propcat  = @("debt","debt","debt","debt","debt","debt","debt")
$category = @("normal","normal","normal","normal","normal","normal","normal")
$date     = @("2011-12-26","2011-12-26","2011-12-26","2011-12-26","2011-12-26","2011-12-26","2011-12-26")
$legname  = @("黃昭順","黃昭順","黃昭順","黃昭順","黃昭順","黃昭順","黃昭順")
$legid    = @(665,665,665,665,665,665,665)
$srcfile  = @("tmp43441","tmp43441","tmp43441","tmp43441","tmp43441","tmp43441","tmp43441")
$idx      = @(116,117,118,119,120,121,122)

for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 2).Value = $species[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 3).Value = $debtor[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 4).Value = $owner[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 5).Value = $total[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 6).Value = $regdate[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 7).Value = $regreas[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 8).Value = $propcat[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 9).Value = $category[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 10).Value = $date[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 11).Value = $legname[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 12).Value = $legid[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 13).Value = $srcfile[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 14).Value = $idx[$i] }
for ($i = 0; $i -lt 7; $i++) { $ws.Cells.Item($i + 2, 1).Value = $idx[$i] }

# --- Clone cell styles from the "保險" sheet so styles.xml stays untouched ---
# Header row style (bold, centered, bordered) -> B1:N1
$template.Range("B1").Copy()
$ws.Range("B1:N1").PasteSpecial(-4122)

# Index-column style (A column data cells) -> A2:A8
$template.Range("A2").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)

# Regular data-cell style -> B2:N8 (also normalises J back off the "@" format)
$template.Range("B2").Copy()
$ws.Range("B2:N8").PasteSpecial(-4122)
